# Auto-generated edit script: updates market price data cells per diff
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H21").Value = 42309.31
$ws.Range("J21").Value = 41668.418
$ws.Range("L21").Value = 41668.418
$ws.Range("N21").Value = -42604.418
$ws.Range("H23").Value = 42309.31
$ws.Range("J23").Value = 41668.418
$ws.Range("L23").Value = 41668.418
$ws.Range("N23").Value = -42136.418
$ws.Range("H29").Value = 5
$ws.Range("I29").Value = 5
$ws.Range("K29").Value = 15
$ws.Range("M29").Value = 266
$ws.Range("H38").Value = 81.375
$ws.Range("I38").Value = 81.375
$ws.Range("J38").Value = 0
$ws.Range("K38").Value = 244.125
$ws.Range("L38").Value = 0
$ws.Range("M38").Value = 127.875
$ws.Range("N38").ClearContents()
$ws.Range("H58").Value = 2340.5454
$ws.Range("J58").Value = 2404.3809
$ws.Range("L58").Value = 7213.1427
$ws.Range("N58").Value = -7513.1427
$ws.Range("H87").Value = 26562.5
$ws.Range("J87").Value = 26562.5
$ws.Range("L87").Value = 26562.5
$ws.Range("N87").Value = -29058.5
$ws.Range("H90").Value = 26562.5
$ws.Range("J90").Value = 26562.5
$ws.Range("L90").Value = 79687.5
$ws.Range("N90").Value = -92167.5
$ws.Range("H129").Value = 1126.6719
$ws.Range("J129").Value = 1151.8871
$ws.Range("L129").Value = 3455.6613
$ws.Range("N129").Value = -13455.6613
$ws.Range("H138").Value = 3415.2273
$ws.Range("I138").Value = 2145.8
$ws.Range("J138").Value = 3967.152
$ws.Range("K138").Value = 6437.400000000001
$ws.Range("L138").Value = 11901.456
$ws.Range("M138").Value = -1297.400000000001
$ws.Range("N138").Value = -22181.456

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H63").Value = 4970.5
$ws.Range("I63").Value = 3941
$ws.Range("J63").Value = 6000
$ws.Range("K63").Value = 3941
$ws.Range("L63").Value = 6000
$ws.Range("M63").Value = -3255
$ws.Range("N63").Value = -7372
$ws.Range("H66").Value = 4970.5
$ws.Range("I66").Value = 3941
$ws.Range("J66").Value = 6000
$ws.Range("K66").Value = 19705
$ws.Range("L66").Value = 30000
$ws.Range("M66").Value = -16273
$ws.Range("N66").Value = -36864
$ws.Range("H132").Value = 4251.2046
$ws.Range("I132").Value = 4712.2
$ws.Range("J132").Value = 2458.4443
$ws.Range("K132").Value = 14136.6
$ws.Range("L132").Value = 7375.3329
$ws.Range("M132").Value = -11606.6
$ws.Range("N132").Value = -12435.3329

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 1546.5883
$ws.Range("I99").Value = 1469.2307
$ws.Range("J99").Value = 1798
$ws.Range("K99").Value = 1469.2307
$ws.Range("L99").Value = 1798
$ws.Range("M99").Value = 28.76929999999993
$ws.Range("N99").Value = -4794
$ws.Range("H140").Value = 88740
$ws.Range("J140").Value = 88740
$ws.Range("L140").Value = 88740
$ws.Range("N140").Value = -99100

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2403.0518
$ws.Range("I31").Value = 1959.4286
$ws.Range("J31").Value = 3078.1304
$ws.Range("K31").Value = 1959.4286
$ws.Range("L31").Value = 3078.1304
$ws.Range("M31").Value = -1664.4286
$ws.Range("N31").Value = -3668.1304
$ws.Range("H34").Value = 2403.0518
$ws.Range("I34").Value = 1959.4286
$ws.Range("J34").Value = 3078.1304
$ws.Range("K34").Value = 1959.4286
$ws.Range("L34").Value = 3078.1304
$ws.Range("M34").Value = -1757.4286
$ws.Range("N34").Value = -3482.1304
$ws.Range("H132").Value = 332615.06
$ws.Range("I132").Value = 468566.22
$ws.Range("J132").Value = 4066.5
$ws.Range("K132").Value = 1405698.66
$ws.Range("L132").Value = 12199.5
$ws.Range("M132").Value = -1403168.66
$ws.Range("N132").Value = -17259.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1965.5
$ws.Range("I5").Value = 1965.5
$ws.Range("J5").Value = 0
$ws.Range("K5").Value = 5896.5
$ws.Range("L5").Value = 0
$ws.Range("M5").Value = -5784.5
$ws.Range("N5").ClearContents()
$ws.Range("H133").Value = 6498.75
$ws.Range("J133").Value = 7993.1665
$ws.Range("L133").Value = 23979.4995
$ws.Range("N133").Value = -34099.49950000001
$ws.Range("H134").Value = 3345.3
$ws.Range("I134").Value = 2208.524
$ws.Range("J134").Value = 5997.778
$ws.Range("K134").Value = 6625.572
$ws.Range("L134").Value = 17993.334
$ws.Range("M134").Value = -1555.572
$ws.Range("N134").Value = -28133.334
$ws.Range("H135").Value = 1965.5
$ws.Range("I135").Value = 1965.5
$ws.Range("J135").Value = 0
$ws.Range("K135").Value = 17689.5
$ws.Range("L135").Value = 0
$ws.Range("M135").Value = -15154.5
$ws.Range("N135").ClearContents()
$ws.Range("H137").Value = 47622988
$ws.Range("I137").Value = 2843.3333
$ws.Range("J137").Value = 83338100
$ws.Range("K137").Value = 8529.999899999999
$ws.Range("L137").Value = 250014300
$ws.Range("M137").Value = -3429.999899999999
$ws.Range("N137").Value = -250024500
$ws.Range("H139").Value = 1799.5
$ws.Range("I139").Value = 1457.7142
$ws.Range("J139").Value = 2995.75
$ws.Range("K139").Value = 4373.142599999999
$ws.Range("L139").Value = 8987.25
$ws.Range("M139").Value = 766.8574000000008
$ws.Range("N139").Value = -19267.25
$ws.Range("H140").Value = 1998.3334
$ws.Range("I140").Value = 998
$ws.Range("J140").Value = 3248.75
$ws.Range("K140").Value = 2994
$ws.Range("L140").Value = 9746.25
$ws.Range("M140").Value = 2186
$ws.Range("N140").Value = -20106.25
$ws.Range("H141").Value = 14009
$ws.Range("I141").Value = 4507.5
$ws.Range("J141").Value = 26677.666
$ws.Range("K141").Value = 13522.5
$ws.Range("L141").Value = 80032.99800000001
$ws.Range("M141").Value = -8342.5
$ws.Range("N141").Value = -90392.99800000001

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 2569.4736
$ws.Range("I132").Value = 1971.8148
$ws.Range("K132").Value = 5915.4444
$ws.Range("M132").Value = -3385.4444

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H55").Value = 232
$ws.Range("I55").Value = 217.92857
$ws.Range("J55").Value = 249.90909
$ws.Range("K55").Value = 217.92857
$ws.Range("L55").Value = 249.90909
$ws.Range("M55").Value = -44.92857000000001
$ws.Range("N55").Value = -595.90909
$ws.Range("H132").Value = 5555.72
$ws.Range("I132").Value = 5883.0586
$ws.Range("J132").Value = 4860.125
$ws.Range("K132").Value = 17649.1758
$ws.Range("L132").Value = 14580.375
$ws.Range("M132").Value = -15119.1758
$ws.Range("N132").Value = -19640.375
$ws.Range("H136").Value = 2559.1333
$ws.Range("I136").Value = 1435.875
$ws.Range("J136").Value = 3842.8572
$ws.Range("K136").Value = 4307.625
$ws.Range("L136").Value = 11528.5716
$ws.Range("M136").Value = -1757.625
$ws.Range("N136").Value = -16628.5716

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 2904.5715
$ws.Range("I136").Value = 2632.5833
$ws.Range("J136").Value = 3267.2222
$ws.Range("K136").Value = 7897.749899999999
$ws.Range("L136").Value = 9801.6666
$ws.Range("M136").Value = -5347.749899999999
$ws.Range("N136").Value = -14901.6666
